# The deck currently applies the "Integral" theme (green palette) to every
# slide via the (single) Slide Master / theme part. The commit swaps the
# palette that is actually applied to the deck back to the default
# "Office Theme" palette (the blue/orange palette that, before this edit,
# only sat - unused - in the Notes Master's theme part).
#
# Re-create that effect by rewriting the 12 theme colors (the DrawingML
# clrScheme slots: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) that are
# exposed through the PowerPoint object model as Slide.ThemeColorScheme.
# RGB values below are standard VBA RGB() long values (0xBBGGRR).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
